$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1939383.2
$ws.Range("J17").Value = 2004022.8
$ws.Range("L17").Value = 6012068.4
$ws.Range("N17").Value = -6012404.4
$ws.Range("H32").Value = 6053.2144
$ws.Range("I32").Value = 5898.8
$ws.Range("J32").Value = 6139
$ws.Range("K32").Value = 5898.8
$ws.Range("L32").Value = 6139
$ws.Range("M32").Value = -5572.8
$ws.Range("N32").Value = -6791
$ws.Range("H64").Value = 4482.25
$ws.Range("J64").Value = 4876.778
$ws.Range("L64").Value = 4876.778
$ws.Range("N64").Value = -5372.778
$ws.Range("H67").Value = 4482.25
$ws.Range("J67").Value = 4876.778
$ws.Range("L67").Value = 4876.778
$ws.Range("N67").Value = -6592.778
$ws.Range("H74").Value = 4713.4546
$ws.Range("I74").Value = 4199.7144
$ws.Range("K74").Value = 4199.7144
$ws.Range("M74").Value = -3263.7144
$ws.Range("H77").Value = 4713.4546
$ws.Range("I77").Value = 4199.7144
$ws.Range("K77").Value = 20998.572
$ws.Range("M77").Value = -16318.572
$ws.Range("H92").Value = 1757.2
$ws.Range("I92").Value = 1420.1428
$ws.Range("J92").Value = 2543.6667
$ws.Range("K92").Value = 1420.1428
$ws.Range("L92").Value = 2543.6667
$ws.Range("M92").Value = -172.1428000000001
$ws.Range("N92").Value = -5039.6667
$ws.Range("H112").Value = 2047.8286
$ws.Range("J112").Value = 3149.3333
$ws.Range("L112").Value = 9447.999899999999
$ws.Range("N112").Value = -11663.9999
$ws.Range("H113").Value = 8091.1333
$ws.Range("I113").Value = 11608
$ws.Range("J113").Value = 2815.8333
$ws.Range("K113").Value = 11608
$ws.Range("L113").Value = 2815.8333
$ws.Range("M113").Value = -8354
$ws.Range("N113").Value = -9323.8333
$ws.Range("H125").Value = 7695.385
$ws.Range("I125").Value = 7168.5713
$ws.Range("J125").Value = 8310
$ws.Range("K125").Value = 64517.14169999999
$ws.Range("L125").Value = 74790
$ws.Range("M125").Value = -62057.14169999999
$ws.Range("N125").Value = -79710
$ws.Range("H127").Value = 455641.2
$ws.Range("I127").Value = 556476.25
$ws.Range("J127").Value = 1883.5
$ws.Range("K127").Value = 1669428.75
$ws.Range("L127").Value = 5650.5
$ws.Range("M127").Value = -1664468.75
$ws.Range("N127").Value = -15570.5
$ws.Range("H135").Value = 689.8421
$ws.Range("I135").Value = 672.6111
$ws.Range("K135").Value = 6053.4999
$ws.Range("M135").Value = -3518.4999
$ws.Range("H141").Value = 3326.1428
$ws.Range("I141").Value = 3089.9412
$ws.Range("K141").Value = 9269.8236
$ws.Range("M141").Value = -4089.8236

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3274.4285
$ws.Range("I45").Value = 2313
$ws.Range("J45").Value = 5678
$ws.Range("K45").Value = 2313
$ws.Range("L45").Value = 5678
$ws.Range("M45").Value = -1936
$ws.Range("N45").Value = -6432
$ws.Range("H52").Value = 34999
$ws.Range("J52").Value = 34999
$ws.Range("L52").Value = 34999
$ws.Range("N52").Value = -35635
$ws.Range("H74").Value = 27501.625
$ws.Range("I74").Value = 2190.625
$ws.Range("K74").Value = 2190.625
$ws.Range("M74").Value = -1316.625
$ws.Range("H77").Value = 27501.625
$ws.Range("I77").Value = 2190.625
$ws.Range("K77").Value = 10953.125
$ws.Range("M77").Value = -6585.125
$ws.Range("H102").Value = 3064.0908
$ws.Range("I102").Value = 2715.25
$ws.Range("K102").Value = 2715.25
$ws.Range("M102").Value = -1093.25
$ws.Range("H132").Value = 8833204
$ws.Range("I132").Value = 3871.25
$ws.Range("K132").Value = 11613.75
$ws.Range("M132").Value = -9083.75
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 10089.333
$ws.Range("I22").Value = 10870.182
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 10870.182
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -10697.182
$ws.Range("N22").Value = -1846
$ws.Range("H35").Value = 34689
$ws.Range("I35").Value = 27533.5
$ws.Range("K35").Value = 27533.5
$ws.Range("M35").Value = -27223.5
$ws.Range("H94").Value = 550.6875
$ws.Range("I94").Value = 370.2
$ws.Range("J94").Value = 851.5
$ws.Range("K94").Value = 370.2
$ws.Range("L94").Value = 851.5
$ws.Range("M94").Value = 80.80000000000001
$ws.Range("N94").Value = -1753.5
$ws.Range("H105").Value = 1183
$ws.Range("I105").Value = 1231.4286
$ws.Range("J105").Value = 674.5
$ws.Range("K105").Value = 1231.4286
$ws.Range("L105").Value = 674.5
$ws.Range("M105").Value = 515.5714
$ws.Range("N105").Value = -4168.5
$ws.Range("H134").Value = 186432.72
$ws.Range("I134").Value = 540000
$ws.Range("K134").Value = 1620000
$ws.Range("M134").Value = -1617465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13444.056
$ws.Range("I58").Value = 4249.3076
$ws.Range("K58").Value = 4249.3076
$ws.Range("M58").Value = -4046.3076
$ws.Range("H99").Value = 13298
$ws.Range("I99").Value = 4622.5
$ws.Range("K99").Value = 4622.5
$ws.Range("M99").Value = -3124.5
$ws.Range("H126").Value = 13298
$ws.Range("I126").Value = 4622.5
$ws.Range("K126").Value = 13867.5
$ws.Range("M126").Value = -11397.5
$ws.Range("H134").Value = 27032262
$ws.Range("I134").Value = 1687.0416
$ws.Range("K134").Value = 5061.1248
$ws.Range("M134").Value = -2526.1248
$ws.Range("H136").Value = 13444.056
$ws.Range("I136").Value = 4249.3076
$ws.Range("K136").Value = 12747.9228
$ws.Range("M136").Value = -10197.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3411.55
$ws.Range("I113").Value = 5311.273
$ws.Range("K113").Value = 15933.819
$ws.Range("M113").Value = -13763.819
$ws.Range("H116").Value = 14594780
$ws.Range("I116").Value = 25414242
$ws.Range("J116").Value = 168831.67
$ws.Range("K116").Value = 76242726
$ws.Range("L116").Value = 506495.01
$ws.Range("M116").Value = -76239284
$ws.Range("N116").Value = -513379.01

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 813.26666
$ws.Range("I97").Value = 599.25
$ws.Range("J97").Value = 1669.3334
$ws.Range("K97").Value = 599.25
$ws.Range("L97").Value = 1669.3334
$ws.Range("M97").Value = -103.25
$ws.Range("N97").Value = -2661.3334
$ws.Range("H132").Value = 2437729.2
$ws.Range("I132").Value = 3641.3333
$ws.Range("J132").Value = 4263295
$ws.Range("K132").Value = 10923.9999
$ws.Range("L132").Value = 12789885
$ws.Range("M132").Value = -8393.999899999999
$ws.Range("N132").Value = -12794945

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2799.5386
$ws.Range("J46").Value = 3219.4
$ws.Range("L46").Value = 3219.4
$ws.Range("N46").Value = -3595.4
$ws.Range("H93").Value = 7307.6924
$ws.Range("I93").Value = 7524.75
$ws.Range("J93").Value = 6960.4
$ws.Range("K93").Value = 7524.75
$ws.Range("L93").Value = 6960.4
$ws.Range("M93").Value = -6276.75
$ws.Range("N93").Value = -9456.4
$ws.Range("H100").Value = 3134.1428
$ws.Range("I100").Value = 2719.25
$ws.Range("J100").Value = 3687.3333
$ws.Range("K100").Value = 2719.25
$ws.Range("L100").Value = 3687.3333
$ws.Range("M100").Value = -2178.25
$ws.Range("N100").Value = -4769.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 641997.75
$ws.Range("I132").Value = 1377.4445
$ws.Range("J132").Value = 1362695.6
$ws.Range("K132").Value = 4132.333500000001
$ws.Range("L132").Value = 4088086.8
$ws.Range("M132").Value = -1602.333500000001
$ws.Range("N132").Value = -4093146.8
